$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dealer")

# Row 5: I5 text stays "ID_Name" (no content change, just internal string-table renumbering)

# Row 8: clear the Source-3 columns (G/H/I) -- fewer sources needed
$ws.Range("G8:I8").ClearContents()

# Row 10: clear the Source-3 columns (G/H/I)
$ws.Range("G10:I10").ClearContents()

# Row 16: Source-2 (D/E/F) now holds what used to be in Source-3 (G/H/I);
# Source-3 columns are cleared.
$ws.Range("D16").Value = "SEIS732_Team_02_Sales_Org"
$ws.Range("E16").Value = "Dealer"
$ws.Range("F16").Value = "DLR_Phone"
$ws.Range("G16:I16").ClearContents()

# Update the active selection to match the author's final state
$ws.Range("G8:I8").Select()
